$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = 120
$ws.Range("H9").Value = 30
$ws.Range("H8").Value = $null
$ws.Range("I8").Value = $null

$ws.Range("H8").Select()
